# Auto-generated edit script: updates market-data derived columns (H-N)
# for specific rows across multiple worksheets, per scheduled market data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Cells.Item(40, 8).Value = 4234.5454
$ws.Cells.Item(40, 9).Value = 7020
$ws.Cells.Item(40, 10).Value = 1913.3334
$ws.Cells.Item(40, 11).Value = 7020
$ws.Cells.Item(40, 12).Value = 1913.3334
$ws.Cells.Item(40, 13).Value = -6845
$ws.Cells.Item(40, 14).Value = -2263.3334

# Row 43
$ws.Cells.Item(43, 8).Value = 1973.75
$ws.Cells.Item(43, 9).Value = 2138.0667
$ws.Cells.Item(43, 10).Value = 1480.8
$ws.Cells.Item(43, 11).Value = 2138.0667
$ws.Cells.Item(43, 12).Value = 1480.8
$ws.Cells.Item(43, 13).Value = -2069.0667
$ws.Cells.Item(43, 14).Value = -1618.8

# Row 64
$ws.Cells.Item(64, 8).Value = 61751.59
$ws.Cells.Item(64, 9).Value = 335266.66
$ws.Cells.Item(64, 10).Value = 3141.2144
$ws.Cells.Item(64, 11).Value = 335266.66
$ws.Cells.Item(64, 12).Value = 3141.2144
$ws.Cells.Item(64, 13).Value = -335018.66
$ws.Cells.Item(64, 14).Value = -3637.2144

# Row 67
$ws.Cells.Item(67, 8).Value = 61751.59
$ws.Cells.Item(67, 9).Value = 335266.66
$ws.Cells.Item(67, 10).Value = 3141.2144
$ws.Cells.Item(67, 11).Value = 335266.66
$ws.Cells.Item(67, 12).Value = 3141.2144
$ws.Cells.Item(67, 13).Value = -334408.66
$ws.Cells.Item(67, 14).Value = -4857.2144

# Row 74
$ws.Cells.Item(74, 8).Value = 3459.8
$ws.Cells.Item(74, 9).Value = 3199.75
$ws.Cells.Item(74, 10).Value = 4500
$ws.Cells.Item(74, 11).Value = 3199.75
$ws.Cells.Item(74, 12).Value = 4500
$ws.Cells.Item(74, 13).Value = -2263.75
$ws.Cells.Item(74, 14).Value = -6372

# Row 77
$ws.Cells.Item(77, 8).Value = 3459.8
$ws.Cells.Item(77, 9).Value = 3199.75
$ws.Cells.Item(77, 10).Value = 4500
$ws.Cells.Item(77, 11).Value = 15998.75
$ws.Cells.Item(77, 12).Value = 22500
$ws.Cells.Item(77, 13).Value = -11318.75
$ws.Cells.Item(77, 14).Value = -31860

# Row 125
$ws.Cells.Item(125, 8).Value = 5386.154
$ws.Cells.Item(125, 9).Value = 2854
$ws.Cells.Item(125, 10).Value = 7556.5713
$ws.Cells.Item(125, 11).Value = 25686
$ws.Cells.Item(125, 12).Value = 68009.14169999999
$ws.Cells.Item(125, 13).Value = -23226
$ws.Cells.Item(125, 14).Value = -72929.14169999999

# Row 132
$ws.Cells.Item(132, 8).Value = 4171825.5
$ws.Cells.Item(132, 9).Value = 4634834
$ws.Cells.Item(132, 10).Value = 4747.5
$ws.Cells.Item(132, 11).Value = 13904502
$ws.Cells.Item(132, 12).Value = 14242.5
$ws.Cells.Item(132, 13).Value = -13901972
$ws.Cells.Item(132, 14).Value = -19302.5

$ws = $wb.Worksheets.Item("ARM")
# Row 102
$ws.Cells.Item(102, 8).Value = 68762.664
$ws.Cells.Item(102, 9).Value = 102059.9
$ws.Cells.Item(102, 11).Value = 102059.9
$ws.Cells.Item(102, 13).Value = -100437.9

# Row 122
$ws.Cells.Item(122, 8).Value = 2890.923
$ws.Cells.Item(122, 9).Value = 2062.2
$ws.Cells.Item(122, 10).Value = 5653.3335
$ws.Cells.Item(122, 11).Value = 6186.599999999999
$ws.Cells.Item(122, 12).Value = 16960.0005
$ws.Cells.Item(122, 13).Value = -3736.599999999999
$ws.Cells.Item(122, 14).Value = -21860.0005

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Cells.Item(99, 8).Value = 2009.9286
$ws.Cells.Item(99, 9).Value = 1140.1428
$ws.Cells.Item(99, 10).Value = 2299.8572
$ws.Cells.Item(99, 11).Value = 1140.1428
$ws.Cells.Item(99, 12).Value = 2299.8572
$ws.Cells.Item(99, 13).Value = 357.8571999999999
$ws.Cells.Item(99, 14).Value = -5295.8572

# Row 105
$ws.Cells.Item(105, 8).Value = 401798.2
$ws.Cells.Item(105, 9).Value = 335660
$ws.Cells.Item(105, 10).Value = 501005.5
$ws.Cells.Item(105, 11).Value = 335660
$ws.Cells.Item(105, 12).Value = 501005.5
$ws.Cells.Item(105, 13).Value = -333913
$ws.Cells.Item(105, 14).Value = -504499.5

# Row 134
$ws.Cells.Item(134, 8).Value = 15340
$ws.Cells.Item(134, 9).Value = 17267.559
$ws.Cells.Item(134, 10).Value = 4417.1665
$ws.Cells.Item(134, 11).Value = 51802.677
$ws.Cells.Item(134, 12).Value = 13251.4995
$ws.Cells.Item(134, 13).Value = -49267.677
$ws.Cells.Item(134, 14).Value = -18321.4995

$ws = $wb.Worksheets.Item("CRP")
# Row 23
$ws.Cells.Item(23, 8).Value = 400
$ws.Cells.Item(23, 9).Value = 400
$ws.Cells.Item(23, 11).Value = 400
$ws.Cells.Item(23, 13).Value = -160

# Row 27
$ws.Cells.Item(27, 8).Value = 400
$ws.Cells.Item(27, 9).Value = 400
$ws.Cells.Item(27, 11).Value = 400
$ws.Cells.Item(27, 13).Value = -208

# Row 122
$ws.Cells.Item(122, 8).Value = 992.4286
$ws.Cells.Item(122, 9).Value = 970.5789
$ws.Cells.Item(122, 11).Value = 2911.7367
$ws.Cells.Item(122, 13).Value = -461.7366999999999

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Cells.Item(131, 8).Value = 10355.25
$ws.Cells.Item(131, 10).Value = 11489.459
$ws.Cells.Item(131, 12).Value = 34468.377
$ws.Cells.Item(131, 14).Value = -44548.377

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Cells.Item(102, 8).Value = 3299.5454
$ws.Cells.Item(102, 9).Value = 1977.1538
$ws.Cells.Item(102, 10).Value = 5209.6665
$ws.Cells.Item(102, 11).Value = 1977.1538
$ws.Cells.Item(102, 12).Value = 5209.6665
$ws.Cells.Item(102, 13).Value = -355.1538
$ws.Cells.Item(102, 14).Value = -8453.666499999999

$ws = $wb.Worksheets.Item("LTW")
# Row 43
$ws.Cells.Item(43, 8).Value = 737282.9
$ws.Cells.Item(43, 9).Value = 5000000
$ws.Cells.Item(43, 10).Value = 26830
$ws.Cells.Item(43, 11).Value = 5000000
$ws.Cells.Item(43, 12).Value = 26830
$ws.Cells.Item(43, 13).Value = -4999807
$ws.Cells.Item(43, 14).Value = -27216

# Row 55
$ws.Cells.Item(55, 8).Value = 637.4
$ws.Cells.Item(55, 9).Value = 741.96155
$ws.Cells.Item(55, 10).Value = 494.3158
$ws.Cells.Item(55, 11).Value = 741.96155
$ws.Cells.Item(55, 12).Value = 494.3158
$ws.Cells.Item(55, 13).Value = -568.96155
$ws.Cells.Item(55, 14).Value = -840.3158000000001

# Row 61
$ws.Cells.Item(61, 8).Value = 2179.2856
$ws.Cells.Item(61, 9).Value = 2137.2727
$ws.Cells.Item(61, 10).Value = 2333.3333
$ws.Cells.Item(61, 11).Value = 2137.2727
$ws.Cells.Item(61, 12).Value = 2333.3333
$ws.Cells.Item(61, 13).Value = -1935.2727
$ws.Cells.Item(61, 14).Value = -2737.3333

# Row 68
$ws.Cells.Item(68, 8).Value = 2686.6538
$ws.Cells.Item(68, 9).Value = 1831.25
$ws.Cells.Item(68, 10).Value = 4055.3
$ws.Cells.Item(68, 11).Value = 1831.25
$ws.Cells.Item(68, 12).Value = 4055.3
$ws.Cells.Item(68, 13).Value = -1082.25
$ws.Cells.Item(68, 14).Value = -5553.3

# Row 71
$ws.Cells.Item(71, 8).Value = 2686.6538
$ws.Cells.Item(71, 9).Value = 1831.25
$ws.Cells.Item(71, 10).Value = 4055.3
$ws.Cells.Item(71, 11).Value = 9156.25
$ws.Cells.Item(71, 12).Value = 20276.5
$ws.Cells.Item(71, 13).Value = -5412.25
$ws.Cells.Item(71, 14).Value = -27764.5

# Row 82
$ws.Cells.Item(82, 8).Value = 2049.2778
$ws.Cells.Item(82, 9).Value = 1580.2
$ws.Cells.Item(82, 10).Value = 2635.625
$ws.Cells.Item(82, 11).Value = 1580.2
$ws.Cells.Item(82, 12).Value = 2635.625
$ws.Cells.Item(82, 13).Value = -1219.2
$ws.Cells.Item(82, 14).Value = -3357.625

# Row 85
$ws.Cells.Item(85, 8).Value = 2049.2778
$ws.Cells.Item(85, 9).Value = 1580.2
$ws.Cells.Item(85, 10).Value = 2635.625
$ws.Cells.Item(85, 11).Value = 1580.2
$ws.Cells.Item(85, 12).Value = 2635.625
$ws.Cells.Item(85, 13).Value = -332.2
$ws.Cells.Item(85, 14).Value = -5131.625

# Row 100
$ws.Cells.Item(100, 8).Value = 1927.2142
$ws.Cells.Item(100, 9).Value = 1771.8572
$ws.Cells.Item(100, 10).Value = 2082.5715
$ws.Cells.Item(100, 11).Value = 1771.8572
$ws.Cells.Item(100, 12).Value = 2082.5715
$ws.Cells.Item(100, 13).Value = -1230.8572
$ws.Cells.Item(100, 14).Value = -3164.5715

# Row 113
$ws.Cells.Item(113, 8).Value = 2179.2856
$ws.Cells.Item(113, 9).Value = 2137.2727
$ws.Cells.Item(113, 10).Value = 2333.3333
$ws.Cells.Item(113, 11).Value = 2137.2727
$ws.Cells.Item(113, 12).Value = 2333.3333
$ws.Cells.Item(113, 13).Value = 32.72730000000001
$ws.Cells.Item(113, 14).Value = -6673.3333

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Cells.Item(126, 8).Value = 1330.8182
$ws.Cells.Item(126, 9).Value = 1309.3684
$ws.Cells.Item(126, 10).Value = 1466.6666
$ws.Cells.Item(126, 11).Value = 3928.1052
$ws.Cells.Item(126, 12).Value = 4399.9998
$ws.Cells.Item(126, 13).Value = -1458.1052
$ws.Cells.Item(126, 14).Value = -9339.9998
